# 2. kodutöö esimene commit
# Adds a new "Nädal 3" worksheet (copied from "Nädal 2"), fills in week-3
# data for the first two log rows, and clears the remaining placeholder
# rows so they stay blank (but keep their original formatting).

$wb = $excel.ActiveWorkbook

# --- 1. Create the new sheet as a copy of "Nädal 2", placed right after it ---
$week2 = $wb.Worksheets.Item("Nädal 2")
$week2.Copy($null, $week2)
$week3 = $wb.Worksheets.Item($wb.Worksheets.Count)
$week3.Name = "Nädal 3"

# --- 2. Clear out the sample/placeholder rows (8-13) on the new sheet,    ---
#        keeping cell formatting intact - only row 6 and 7 keep real data.
$week3.Range("B8:J13").ClearContents()

# --- 3. Row 6: first log entry of week 3 ---
# Column B keeps a literal text date (not a real Excel date) just like the
# source row did, so force Text format before typing it in.
$week3.Range("B6").NumberFormat = "@"
$week3.Range("B6").Value = "11.02.2020"
$week3.Range("B6").NumberFormat = "dd/mm/yyyy;@"

$week3.Range("D6").Value = 0.39930555555555558
$week3.Range("F6").Value = 95
$week3.Range("H6").Value = "Milleks on õppejõud? "

# --- 4. Row 7: second log entry of week 3 ---
$week3.Range("B7").Value = 43874
$week3.Range("C7").Value = 0.76041666666666663
$week3.Range("D7").Value = 0.80208333333333337
$week3.Range("E7").ClearContents()
$week3.Range("F7").Value = 80
$week3.Range("G7").Value = "PodCast"
$week3.Range("H7").Value = "Algorütm 2"

# --- 5. View state: leave "Nädal 2" with its whole table selected, then  ---
#        land on "Nädal 3" with H7 selected as the new active sheet.
$week2.Activate()
$week2.Range("A1:J18").Select()

$week3.Activate()
$week3.Range("H7").Select()
